$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers
$ws.Range("A1").Value = "Job_Id"
$ws.Range("B1").Value = "Job_Title"
$ws.Range("C1").Value = "Job_Description"
$ws.Range("D1").Value = "Total_Years_Min_Exp"
$ws.Range("E1").Value = "Total_Years_Max_Exp"
$ws.Range("F1").Value = "LinkedIn_Poster"
$ws.Range("G1").Value = "LinkedIn_Posted"
$ws.Range("H1").Value = "Resume_received"
$ws.Range("I1").Value = "Resume_downloaded"

# Row 2 data
$ws.Range("A2").Value = "JD_001"
$ws.Range("B2").Value = "Cyber Security Engineer"
$ws.Range("C2").Value = "We are seeking a Cyber Security Engineer to build and maintain high-quality software solutions.`nWork with global teams to drive innovation and deliver scalable applications.`nJoin Akkodis and be part of a tech-driven, collaborative environment."
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 4

# Restore default (non-custom) row height after the multi-line text auto-expanded it
$ws.Rows.Item(2).AutoFit()
